# Adds the Pearson correlation coefficient output (cor.test) plus the
# interpreting paragraph right after the "... no hay un correlación entre
# las dos variables." paragraph, before the "Construya un gráfico ..."
# list item.

$d = $word.ActiveDocument

$anchorText = "Construya un gráfico de dispersión que relacione la variable eruption"
$sourceCodeText = "## [1] 0.9008112"
$interpretationText = "A través del indice de correlación obtenido (0,9) se puede considerar que existe una correlación positiva entre las variables."

# --- New "SourceCode" paragraph with the cor.test() console output -------
$rng = $d.Content
$rng.Find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Collapse(1)
$rng.InsertParagraphBefore()
$rng.Style = "SourceCode"
$rng.Text = $sourceCodeText
$rng.Style = "VerbatimChar"

# --- New "FirstParagraph" paragraph interpreting the correlation index ---
$rng2 = $d.Content
$rng2.Find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng2.Collapse(1)
$rng2.InsertParagraphBefore()
$rng2.Style = "FirstParagraph"
$rng2.Text = $interpretationText

Write-Output "Inserted correlation coefficient paragraphs."
